$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3033
$ws1.Range("F5").Value = 420
$ws1.Range("F6").Value = 285

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3033
$ws4.Range("F5").Value = 420
$ws4.Range("F7").Value = 285
